# Revert "adding term 2.0.0"
#
# 1. Restore the Metadata sheet's Version/Date/Contact values to the
#    pre-2.0.0 release (1.1.0).
# 2. Restore the "descendent-of" concept value on the "Include from FSIII"
#    sheet back to "B".
# 3. Remove the duplicated "Include from FSIII 2" worksheet that was added
#    alongside the 2.0.0 term.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B3").Value = "1.1.0"
$metadata.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$metadata.Range("B10").Value = "No display for ContactDetail"

$include1 = $wb.Worksheets.Item("Include from FSIII")
$include1.Range("C2").Value = "B"

$include2 = $wb.Worksheets.Item("Include from FSIII 2")
$include2.Delete()
